$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (374) down into the new rows (375-385)
$ws.Range("A374:D374").Copy() | Out-Null
$ws.Range("A375:D385").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(375, 44449, 1, 2, 12.52348152786475),
    @(376, 44450, 2, 4, 25.04696305572949),
    @(377, 44451, 0, 4, 25.04696305572949),
    @(378, 44452, 1, 4, 25.04696305572949),
    @(379, 44453, 0, 4, 25.04696305572949),
    @(380, 44454, 0, 4, 25.04696305572949),
    @(381, 44455, 0, 4, 25.04696305572949),
    @(382, 44456, 0, 3, 18.78522229179712),
    @(383, 44457, 1, 2, 12.52348152786475),
    @(384, 44458, 1, 3, 18.78522229179712),
    @(385, 44459, 0, 2, 12.52348152786475)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
